$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TheWayWeLiveNowV2")
$ws1.Columns.Item(1).Hidden = $false
for ($i=0; $i -lt 30; $i++) {
    $test = 16.05 + $i*0.005
    $ws1.Columns.Item(1).ColumnWidth = $test
    Write-Output "$test -> $($ws1.Columns.Item(1).ColumnWidth)"
}
